# RedemptionsReporting.xlsx — roll the report from store 98 (Waterlooville)
# / period 01-Jul-2024..07-Jul-2024 to store 388 (Lancashire Golf Academy)
# / period 03-Jun-2024..09-Jun-2024, and record the resent statement
# (Statement Date / Statement No. / Statement Amt) for the redemption row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block --------------------------------------------------------
$ws.Range("B2").Value = "388 | Lancashire Golf Academy"
$ws.Range("B4").Value = "03-Jun-2024 To 09-Jun-2024"

# --- Summary figures (rows 5 & 6) ----------------------------------------
$ws.Range("B5").Value = -75
$ws.Range("B6").Value = -63.3

# --- Detail row 8 ----------------------------------------------------------
$ws.Range("A8").Value = 45446.84405092592
$ws.Range("B8").Value = -75
$ws.Range("D8").Value = -9.75
$ws.Range("E8").Value = -1.95
$ws.Range("F8").Value = -63.3
$ws.Range("G8").Value = 388
$ws.Range("H8").Value = "Lancashire Golf Academy"
$ws.Range("I8").Value = "BB1 9LF"
$ws.Range("J8").Value = "EP0130108000000525"

# K8 ("Card ID") is a numeric-looking value that must stay text. Writing it
# with a leading apostrophe forces text storage; re-apply the cell's own
# original (non quote-prefixed) style afterwards so formatting doesn't drift.
$k8Style = $ws.Range("K8").Style
$ws.Range("K8").Value = "'1241458099"
$ws.Range("K8").Style = $k8Style

# Statement date / number / amount for the resent statement (feedback resend + bcc)
$ws.Range("L8").Value = 45491.647632025466
$ws.Range("M8").Value = "GCP0000060"
$ws.Range("N8").Value = 63.3

# --- Totals row 9 ----------------------------------------------------------
$ws.Range("B9").Value = -75
$ws.Range("D9").Value = -9.75
$ws.Range("E9").Value = -1.95
$ws.Range("F9").Value = -63.3

# --- Column width tweaks ---------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 30.3
$ws.Columns.Item(8).ColumnWidth = 21.92
